$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the B-column "current" values on Login Page to what used to be the
# (correct) C-column values, then clear out the now-redundant C column.
# (Using .Formula rather than .Value keeps the existing cell style index
# instead of Excel minting new duplicate cellXfs entries for the re-applied
# Hyperlink style.)
$ws.Range("B2").Formula = "wipro@mailinator.com"
$ws.Range("B4").Formula = "automate@workstreets.com"
$ws.Range("B6").Formula = "source1@mailinator.com"
$ws.Range("B8").Formula = "barish@nada.email"

$ws.Range("C2:C9").ClearContents()

# Hyperlinks.Delete() on this runtime drops the whole sheet's hyperlink
# collection no matter which range it is invoked from, so rebuild only the
# links that should remain (column B except B6) with correct targets.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:wipro@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:automate@workstreets.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:barish@nada.email") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:admin@123") | Out-Null

# Match the new selection left behind on the sheet (user selected C2:C11
# before clearing it).
$ws.Range("C2:C11").Select()
